$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999147019625,
    0.9971433942856277,
    0.98549050094014135,
    0.98503563183308196,
    0.97131309834847968,
    0.94974408631483009,
    0.94475851788178089,
    0.93820552771491283,
    0.93080961528061945,
    0.92433020953261213,
    0.92342945941248311,
    0.92197138653539967,
    0.91690235917650509,
    0.91546606631013461,
    0.91545116901304624,
    0.91294478377798938,
    0.90923698757598559,
    0.90812809508204406,
    0.9960902384038528,
    0.98897339181316957,
    0.98757493608436109,
    0.98631044112701072,
    0.97845582473702197,
    0.96543506421546965,
    0.95897809592949346,
    0.94799643244531917,
    0.94581030311236391,
    0.9368514875180195,
    0.93088945587174909,
    0.92890631350126895,
    0.93316235733823893,
    0.93564539215791653,
    0.93512539668193218
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
